$d = $word.ActiveDocument

# Remove the "Предметни наставник:" text that used to sit in the first
# table cell (the cell's paragraph/formatting is left in place - only the
# run carrying the text goes away).
$d.Content.Find.Execute("Предметни наставник:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
